$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free inline writes below. For values that look numeric but must stay
# as text (leading zeros), we format the cell as Text, assign, then clear the
# formatting again so no stray number format lingers on the cell itself.

# --- Row 8: School ID 2 -> Global Reciprocal College (grc admin) ---
$ws.Range("A8").Value = 2
$ws.Range("B8").Value = "Global Reciprocal College"
$ws.Range("C8").Value = "public"
$ws.Range("D8").Value = "DCS Valenzuela"
$ws.Range("E8").Value = "Congressional I"
$ws.Range("F8").Value = "grc"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "08745120956230"
$ws.Range("G8").ClearFormats()
$ws.Range("H8").Value = "grc@gmail.com"
$ws.Range("I8").Value = "2024-10-24 08:26:33"

# --- Row 9: School ID 3 -> Global Reciprocal College, second submission ---
$ws.Range("A9").Value = 3
$ws.Range("B9").Value = "Global Reciprocal College"
$ws.Range("C9").Value = "public"
$ws.Range("D9").Value = "DCS Valenzuela"
$ws.Range("E9").Value = "Congressional I"
$ws.Range("F9").Value = "grc"
$ws.Range("G9").Value = 9123456893
$ws.Range("H9").Value = "grc@gmail.com	"
$ws.Range("I9").Value = "2024-10-24 10:22:51"

# --- Row 10: School ID 9 -> ELI School (private), no contact info on file ---
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "ELI School"
$ws.Range("C10").Value = "private"
$ws.Range("D10").Value = "DCS Valenzuela"
$ws.Range("E10").Value = "Congressional I"
$ws.Range("F10").Value = ""
$ws.Range("G10").Value = ""
$ws.Range("H10").Value = ""
$ws.Range("I10").Value = "2024-10-25 10:46:59"

# --- Row 11: School ID 555666 -> Maysan National High School ---
$ws.Range("A11").Value = 555666
$ws.Range("B11").Value = "Maysan National High School"
$ws.Range("C11").Value = "public"
$ws.Range("D11").Value = "DCS Valenzuela"
$ws.Range("E11").Value = "Congressional I"
$ws.Range("F11").Value = "Lorem Ipsum A"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "09060158736"
$ws.Range("G11").ClearFormats()
$ws.Range("H11").Value = "exonatural321@gmail.com"
$ws.Range("I11").Value = "2024-09-04 01:50:22"

# --- Row 12 (School ID 999999 / Example) no longer exists: remove it ---
$ws.Rows("12:12").Delete()
